function Replace-SubtitleText($doc, $old, $new) {
    $rng = $doc.Content
    $rng.Start = 0
    $rng.End = $doc.Content.End
    while ($rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $rng.Text = $new
        $newStart = $rng.End
        $rng = $doc.Content
        $rng.Start = $newStart
    }
}

$d = $word.ActiveDocument

Replace-SubtitleText $d 'The playful mathematicians:' 'Wanahisabati wanaocheza:'
Replace-SubtitleText $d '** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino' '** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino'
Replace-SubtitleText $d '[Music]' '[Muziki]'
Replace-SubtitleText $d 'there are two mathematicians, let''s call' 'kuna wanahisabati wawili, tupige simu'
Replace-SubtitleText $d 'them Fil and Mike who meet each other' 'Fil na Mike wanaokutana'
Replace-SubtitleText $d 'again after a long time. After some' 'tena baada ya muda mrefu. Baada ya baadhi'
Replace-SubtitleText $d 'chatting, Phil says he has three children, then' 'kuzungumza, Phil anasema ana watoto watatu, basi'
Replace-SubtitleText $d 'Mike, astonished, asks: ''How old are they?'' Fil,' 'Kwa mshangao, Mike anauliza: ''Wana umri gani?'' Fil,'
Replace-SubtitleText $d 'being a playful mathematician, answers' 'kuwa mwanahisabati mchezaji, anajibu'
Replace-SubtitleText $d '''You tell me! I''ll give you a hint: if you' '''Wewe niambie! Nitakupa kidokezo: ikiwa wewe'
Replace-SubtitleText $d 'multiply the three ages together you' 'zidisheni enzi tatu pamoja ninyi'
Replace-SubtitleText $d 'get 36.'' Mike takes sometimes to think' 'pata 36.'' Mike huchukua wakati mwingine kufikiria'
Replace-SubtitleText $d 'and says: ''I''m sorry Fil, but I do need' 'na kusema: ''Samahani Fil, lakini nahitaji'
Replace-SubtitleText $d 'another hint. So Fil tells Mike:' 'kidokezo kingine. Kwa hivyo Fil anamwambia Mike:'
Replace-SubtitleText $d '''Yes, sure, here it is: if you had up to' '''Ndiyo, hakika, hapa ni: kama alikuwa na hadi'
Replace-SubtitleText $d 'three ages you get the number of math' 'miaka mitatu unapata idadi ya hesabu'
Replace-SubtitleText $d 'papers we publish together. Do you remember it?''' 'karatasi tunachapisha pamoja. Je, unaikumbuka?'''
Replace-SubtitleText $d '''Yes I do remember How many, but still' '''Ndio nakumbuka wangapi, lakini bado'
Replace-SubtitleText $d 'I do not have enough information! I need' 'Sina taarifa za kutosha! nahitaji'
Replace-SubtitleText $d 'at least one more.'' Fil says: ''Yes don''t' 'angalau moja zaidi.'' Fil anasema: ''Ndiyo usifanye hivyo'
Replace-SubtitleText $d 'worry but this is the last one:' 'wasiwasi lakini hii ni ya mwisho:'
Replace-SubtitleText $d 'The youngest one has blues eyes.'' And' 'Mdogo ana macho ya blues.'' Na'
Replace-SubtitleText $d 'suddenly Mike gets the answer. You' 'ghafla Mike anapata jibu. Wewe'
Replace-SubtitleText $d 'hear the conversation but you don''t know' 'sikia mazungumzo lakini hujui'
Replace-SubtitleText $d 'how many papers they published together.' 'ni karatasi ngapi walichapisha pamoja.'
Replace-SubtitleText $d 'However, you do want to know the ages of' 'Hata hivyo, unataka kujua umri wa'
Replace-SubtitleText $d 'the three children. Can you figure them' 'watoto watatu. Je, unaweza kuwahesabu'
Replace-SubtitleText $d 'out?' 'nje?'
